$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=2; Col=4; Value='67.902.75'},
    @{Row=2; Col=5; Value='  +0.07%  '},
    @{Row=3; Col=4; Value='3.333.91'},
    @{Row=3; Col=5; Value='  +0.81%  '},
    @{Row=4; Col=4; Value='''0.998'},
    @{Row=4; Col=5; Value='  +0.08%  '},
    @{Row=5; Col=4; Value='''583.88'},
    @{Row=5; Col=5; Value='  +0.32%  '},
    @{Row=6; Col=4; Value='''175.06'},
    @{Row=6; Col=5; Value='  -0.44%  '},
    @{Row=7; Col=4; Value='''0.998'},
    @{Row=7; Col=5; Value='  -0.05%  '},
    @{Row=8; Col=4; Value='''0.590'},
    @{Row=8; Col=5; Value='  +1.63%  '},
    @{Row=9; Col=4; Value='3.331.36'},
    @{Row=9; Col=5; Value='  +0.87%  '},
    @{Row=10; Col=4; Value='''0.182'},
    @{Row=10; Col=5; Value='  +4.29%  '},
    @{Row=11; Col=4; Value='''0.580'},
    @{Row=11; Col=5; Value='  +1.31%  '},
    @{Row=12; Col=4; Value='''47.29'},
    @{Row=12; Col=5; Value='  +4.14%  '},
    @{Row=13; Col=4; Value='''0.0000273'},
    @{Row=13; Col=5; Value='  +1.63%  '},
    @{Row=14; Col=4; Value='''699.54'},
    @{Row=14; Col=5; Value='  +5.01%  '},
    @{Row=15; Col=4; Value='3.867.63'},
    @{Row=15; Col=5; Value='  +1.20%  '},
    @{Row=16; Col=4; Value='''8.38'},
    @{Row=16; Col=5; Value='  +0.59%  '},
    @{Row=17; Col=4; Value='67.891.68'},
    @{Row=17; Col=5; Value='  +0.19%  '},
    @{Row=18; Col=5; Value='  +0.65%  '},
    @{Row=19; Col=4; Value='3.326.92'},
    @{Row=19; Col=5; Value='  +1.13%  '},
    @{Row=20; Col=4; Value='''17.49'},
    @{Row=20; Col=5; Value='  +0.59%  '},
    @{Row=21; Col=4; Value='''11.13'},
    @{Row=21; Col=5; Value='  +2.77%  '},
    @{Row=22; Col=4; Value='''0.891'},
    @{Row=22; Col=5; Value='  +0.72%  '},
    @{Row=23; Col=4; Value='''5.42'},
    @{Row=23; Col=5; Value='  -0.12%  '},
    @{Row=24; Col=4; Value='''16.93'},
    @{Row=24; Col=5; Value='  -0.90%  '},
    @{Row=25; Col=4; Value='''101.34'},
    @{Row=25; Col=5; Value='  +3.06%  '},
    @{Row=26; Col=5; Value='  +1.28%  '},
    @{Row=27; Col=4; Value='''2.69'},
    @{Row=27; Col=5; Value='  +1.20%  '},
    @{Row=28; Col=4; Value='''9.45'},
    @{Row=28; Col=5; Value='  +2.79%  '},
    @{Row=29; Col=4; Value='''32.93'},
    @{Row=29; Col=5; Value='  +0.45%  '},
    @{Row=30; Col=4; Value='''8.54'},
    @{Row=30; Col=5; Value='  +2.10%  '},
    @{Row=31; Col=4; Value='''6.98'},
    @{Row=31; Col=5; Value='  -0.81%  '},
    @{Row=32; Col=4; Value='''574.15'},
    @{Row=32; Col=5; Value='  -1.31%  '},
    @{Row=33; Col=4; Value='''11.01'},
    @{Row=33; Col=5; Value='  +0.75%  '},
    @{Row=34; Col=5; Value='  +2.09%  '},
    @{Row=35; Col=4; Value='3.731.54'},
    @{Row=35; Col=5; Value='  -0.65%  '},
    @{Row=36; Col=5; Value='  +0.12%  '},
    @{Row=37; Col=4; Value='''56.66'},
    @{Row=37; Col=5; Value='  +2.00%  '},
    @{Row=38; Col=4; Value='''3.28'},
    @{Row=38; Col=5; Value='  -2.94%  '},
    @{Row=39; Col=4; Value='''35.67'},
    @{Row=39; Col=5; Value='  +10.19%  '},
    @{Row=40; Col=5; Value='  +2.30%  '},
    @{Row=41; Col=2; Value='Fetch.AI'},
    @{Row=41; Col=3; Value='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'},
    @{Row=41; Col=4; Value='''2.61'},
    @{Row=41; Col=5; Value='  -0.90%  '},
    @{Row=42; Col=2; Value='Stacks'},
    @{Row=42; Col=3; Value='https://coinranking.com/coin/mMPrMcB7+stacks-stx'},
    @{Row=42; Col=4; Value='''3.14'},
    @{Row=42; Col=5; Value='  +2.99%  '},
    @{Row=43; Col=4; Value='0.0₃0673'},
    @{Row=43; Col=5; Value='  +1.57%  '},
    @{Row=44; Col=4; Value='''0.334'},
    @{Row=44; Col=5; Value='  +1.91%  '},
    @{Row=45; Col=5; Value='  +1.27%  '},
    @{Row=46; Col=4; Value='''0.0406'},
    @{Row=46; Col=5; Value='  +0.85%  '},
    @{Row=47; Col=4; Value='''2.62'},
    @{Row=47; Col=5; Value='  +1.44%  '},
    @{Row=48; Col=4; Value='''0.129'},
    @{Row=48; Col=5; Value='  +1.48%  '},
    @{Row=49; Col=5; Value='  -0.08%  '},
    @{Row=50; Col=5; Value='  -2.25%  '},
    @{Row=51; Col=4; Value='''130.58'},
    @{Row=51; Col=5; Value='  +1.22%  '}
)

foreach ($item in $changes) {
    $ws.Cells.Item($item.Row, $item.Col).Value = $item.Value
}

Write-Host "Applied $($changes.Count) cell updates"